# Add a new "Future Development" slide before the final "Application links"
# slide (Title and Content layout == layout index 2 / ppLayoutText).

$p = $ppt.ActivePresentation

$oldLastIndex = $p.Slides.Count            # index of the existing "Application links" slide (5)
$newSlide = $p.Slides.Add($oldLastIndex + 1, 2)   # create it at the end (position 6) ...

# --- Title -----------------------------------------------------------
$title = $newSlide.Shapes.Item(1).TextFrame.TextRange
$title.Text = "Future Development"
$title.Font.Bold = $true

# --- Body content ------------------------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Would like to incorporate advance battle mechanics`rAn extended story`rSound effects and music`rImproved scoring`rGeneral style and aesthetic changes"

# Paragraph 1: "Would like to incorporate advance battle mechanics" -> cleared formatting
$r = $body.Characters(1, 50)
$r.Font.Bold = $false
$r.Font.Italic = $false
$r.Font.Shadow = $false
$r.Font.Name = "Slack-Lato"

# Paragraph 2: "A" + "n extended story"
$r = $body.Characters(52, 1)
$r.Font.Name = "Slack-Lato"
$r = $body.Characters(53, 16)
$r.Font.Bold = $false
$r.Font.Italic = $false
$r.Font.Shadow = $false
$r.Font.Name = "Slack-Lato"

# Paragraph 3: "S" + "ound effects and music"
$r = $body.Characters(70, 1)
$r.Font.Name = "Slack-Lato"
$r = $body.Characters(71, 22)
$r.Font.Bold = $false
$r.Font.Italic = $false
$r.Font.Shadow = $false
$r.Font.Name = "Slack-Lato"

# Paragraph 4: "Improved scoring"
$r = $body.Characters(94, 16)
$r.Font.Bold = $false
$r.Font.Italic = $false
$r.Font.Shadow = $false
$r.Font.Name = "Slack-Lato"

# Paragraph 5: "General style and " + "aesthetic changes"
$r = $body.Characters(111, 18)
$r.Font.Name = "Slack-Lato"
$r = $body.Characters(129, 17)
$r.Font.Name = "Slack-Lato"

# Trailing blank paragraph (matches the original author leaving an empty line)
$body.InsertAfter("`r") | Out-Null

# --- Re-order: the new slide belongs right before the old last slide ---
$newSlide.MoveTo($oldLastIndex)
